$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1850220264317181
$ws.Range("C2").Value = 0.5682819383259912
$ws.Range("J2").Value = 0.004405286343612335
$ws.Range("P2").Value = 0.13215859030837
$ws.Range("S2").Value = 0.1101321585903084
$ws.Range("C3").Value = 0.01503759398496241
$ws.Range("J3").Value = 0.03759398496240601
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2330827067669173
$ws.Range("J4").Value = 0.1379310344827586
$ws.Range("P4").Value = 0.5172413793103449
$ws.Range("S4").Value = 0.3448275862068966
$ws.Range("B6").Value = 0.05150214592274678
$ws.Range("D6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.07725321888412018
$ws.Range("J6").Value = 0.2575107296137339
$ws.Range("O6").Value = 0.008583690987124463
$ws.Range("Q6").Value = 0.1502145922746781
$ws.Range("R6").Value = 0.07725321888412018
$ws.Range("S6").Value = 0.369098712446352
$ws.Range("B7").Value = 0.1014492753623188
$ws.Range("D7").Value = 0.01449275362318841
$ws.Range("F7").Value = 0.07729468599033816
$ws.Range("J7").Value = 0.106280193236715
$ws.Range("O7").Value = 0.03381642512077294
$ws.Range("Q7").Value = 0.1207729468599034
$ws.Range("R7").Value = 0.06280193236714976
$ws.Range("S7").Value = 0.4830917874396135
$ws.Range("B8").Value = 0.0640495867768595
$ws.Range("D8").Value = 0.01446280991735537
$ws.Range("F8").Value = 0.08264462809917356
$ws.Range("J8").Value = 0.09917355371900827
$ws.Range("O8").Value = 0.01859504132231405
$ws.Range("Q8").Value = 0.1652892561983471
$ws.Range("R8").Value = 0.121900826446281
$ws.Range("S8").Value = 0.4338842975206612
$ws.Range("B9").Value = 0.08670520231213873
$ws.Range("D9").Value = 0.0115606936416185
$ws.Range("F9").Value = 0.04046242774566474
$ws.Range("J9").Value = 0.09826589595375723
$ws.Range("O9").Value = 0.01734104046242774
$ws.Range("Q9").Value = 0.1734104046242775
$ws.Range("R9").Value = 0.07514450867052024
$ws.Range("S9").Value = 0.4971098265895953
$ws.Range("B10").Value = 0.081145584725537
$ws.Range("D10").Value = 0.01352426412092283
$ws.Range("E10").Value = 0.001591089896579157
$ws.Range("F10").Value = 0.06682577565632458
$ws.Range("J10").Value = 0.1121718377088305
$ws.Range("O10").Value = 0.01352426412092283
$ws.Range("Q10").Value = 0.2171837708830549
$ws.Range("R10").Value = 0.07159904534606205
$ws.Range("S10").Value = 0.4224343675417661
$ws.Range("G11").Value = 0.1335227272727273
$ws.Range("J11").Value = 0.1193181818181818
$ws.Range("K11").Value = 0.2159090909090909
$ws.Range("L11").Value = 0.5198863636363636
$ws.Range("S11").Value = 0.01136363636363636
$ws.Range("G12").Value = 0.6914893617021277
$ws.Range("J12").Value = 0.2393617021276596
$ws.Range("K12").Value = 0.005319148936170213
$ws.Range("L12").Value = 0.01063829787234043
$ws.Range("S12").Value = 0.05319148936170213
$ws.Range("G13").Value = 0.7111111111111111
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01951219512195122
$ws.Range("H15").Value = 0.1707317073170732
$ws.Range("I15").Value = 0.07804878048780488
$ws.Range("J15").Value = 0.3170731707317073
$ws.Range("K15").Value = 0.07317073170731707
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("N15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.06341463414634146
$ws.Range("S15").Value = 0.2634146341463415
$ws.Range("F16").Value = 0.007407407407407408
$ws.Range("H16").Value = 0.1777777777777778
$ws.Range("I16").Value = 0.05185185185185185
$ws.Range("J16").Value = 0.3925925925925926
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.02222222222222222
$ws.Range("O16").Value = 0.08148148148148149
$ws.Range("S16").Value = 0.1555555555555556
$ws.Range("F17").Value = 0.009111617312072893
$ws.Range("H17").Value = 0.1913439635535308
$ws.Range("I17").Value = 0.1002277904328018
$ws.Range("J17").Value = 0.4578587699316629
$ws.Range("K17").Value = 0.07972665148063782
$ws.Range("M17").Value = 0.01366742596810934
$ws.Range("O17").Value = 0.05466970387243736
$ws.Range("S17").Value = 0.09339407744874716
$ws.Range("F18").Value = 0.02631578947368421
$ws.Range("H18").Value = 0.2263157894736842
$ws.Range("I18").Value = 0.05789473684210526
$ws.Range("J18").Value = 0.4105263157894737
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.005263157894736842
$ws.Range("O18").Value = 0.07368421052631578
$ws.Range("S18").Value = 0.1
$ws.Range("F19").Value = 0.02027516292541636
$ws.Range("H19").Value = 0.2179580014482259
$ws.Range("I19").Value = 0.06951484431571325
$ws.Range("J19").Value = 0.3482983345401883
$ws.Range("K19").Value = 0.1383055756698045
$ws.Range("M19").Value = 0.02534395365677046
$ws.Range("O19").Value = 0.05865314989138305
$ws.Range("S19").Value = 0.1216509775524982
